$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 150
$ws.Range("I4").Value = 150
$ws.Range("K4").Value = 150
$ws.Range("M4").Value = -36
$ws.Range("H17").Value = 911140.8
$ws.Range("J17").Value = 911140.8
$ws.Range("L17").Value = 2733422.4
$ws.Range("N17").Value = -2733758.4
$ws.Range("H64").Value = 3570.6667
$ws.Range("I64").Value = 3297.5
$ws.Range("J64").Value = 3707.25
$ws.Range("K64").Value = 3297.5
$ws.Range("L64").Value = 3707.25
$ws.Range("M64").Value = -3049.5
$ws.Range("N64").Value = -4203.25
$ws.Range("H67").Value = 3570.6667
$ws.Range("I67").Value = 3297.5
$ws.Range("J67").Value = 3707.25
$ws.Range("K67").Value = 3297.5
$ws.Range("L67").Value = 3707.25
$ws.Range("M67").Value = -2439.5
$ws.Range("N67").Value = -5423.25
$ws.Range("H132").Value = 1200
$ws.Range("I132").Value = 1200
$ws.Range("K132").Value = 3600
$ws.Range("M132").Value = -1070
$ws.Range("H141").Value = 2668.487
$ws.Range("I141").Value = 1851.6154
$ws.Range("K141").Value = 5554.8462
$ws.Range("M141").Value = -374.8462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1670
$ws.Range("I45").Value = 1670
$ws.Range("K45").Value = 1670
$ws.Range("M45").Value = -1293
$ws.Range("H61").Value = 2617.4167
$ws.Range("I61").Value = 2409.5557
$ws.Range("K61").Value = 2409.5557
$ws.Range("M61").Value = -2197.5557
$ws.Range("H110").Value = 1545.5385
$ws.Range("I110").Value = 1424.3334
$ws.Range("K110").Value = 1424.3334
$ws.Range("M110").Value = 620.6666
$ws.Range("H122").Value = 2493.8948
$ws.Range("I122").Value = 2485.9375
$ws.Range("J122").Value = 2536.3333
$ws.Range("K122").Value = 7457.8125
$ws.Range("L122").Value = 7608.999899999999
$ws.Range("M122").Value = -5007.8125
$ws.Range("N122").Value = -12508.9999
$ws.Range("H132").Value = 6201.7144
$ws.Range("I132").Value = 6195.6665
$ws.Range("K132").Value = 18586.9995
$ws.Range("M132").Value = -16056.9995
$ws.Range("H136").Value = 2617.4167
$ws.Range("I136").Value = 2409.5557
$ws.Range("K136").Value = 7228.6671
$ws.Range("M136").Value = -4678.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1937.5
$ws.Range("I105").Value = 1875
$ws.Range("K105").Value = 1875
$ws.Range("M105").Value = -128
$ws.Range("H134").Value = 6105.273
$ws.Range("I134").Value = 6782.0625
$ws.Range("J134").Value = 4300.5
$ws.Range("K134").Value = 20346.1875
$ws.Range("L134").Value = 12901.5
$ws.Range("M134").Value = -17811.1875
$ws.Range("N134").Value = -17971.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 3500
$ws.Range("I47").Value = 3500
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 3500
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -2934
$ws.Range("N47").ClearContents()
$ws.Range("H58").Value = 2276.5
$ws.Range("I58").Value = 2283.1428
$ws.Range("J58").Value = 2261
$ws.Range("K58").Value = 2283.1428
$ws.Range("L58").Value = 2261
$ws.Range("M58").Value = -2080.1428
$ws.Range("N58").Value = -2667
$ws.Range("H59").Value = 28508.8
$ws.Range("I59").Value = 15104
$ws.Range("K59").Value = 15104
$ws.Range("M59").Value = -13959
$ws.Range("H68").Value = 28331.584
$ws.Range("J68").Value = 29998
$ws.Range("L68").Value = 29998
$ws.Range("N68").Value = -31496
$ws.Range("H71").Value = 28331.584
$ws.Range("J71").Value = 29998
$ws.Range("L71").Value = 89994
$ws.Range("N71").Value = -97482
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H132").Value = 4219.5
$ws.Range("J132").Value = 5036.2
$ws.Range("L132").Value = 15108.6
$ws.Range("N132").Value = -20168.6
$ws.Range("H136").Value = 2276.5
$ws.Range("I136").Value = 2283.1428
$ws.Range("J136").Value = 2261
$ws.Range("K136").Value = 6849.428400000001
$ws.Range("L136").Value = 6783
$ws.Range("M136").Value = -4299.428400000001
$ws.Range("N136").Value = -11883

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 397.83334
$ws.Range("J114").Value = 271.75
$ws.Range("L114").Value = 815.25
$ws.Range("N114").Value = -7323.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 204.66667
$ws.Range("I2").Value = 269.66666
$ws.Range("J2").Value = 139.66667
$ws.Range("K2").Value = 269.66666
$ws.Range("L2").Value = 139.66667
$ws.Range("M2").Value = -156.66666
$ws.Range("N2").Value = -365.66667
$ws.Range("H122").Value = 7759.2
$ws.Range("I122").Value = 8642.75
$ws.Range("K122").Value = 25928.25
$ws.Range("M122").Value = -23478.25
$ws.Range("H132").Value = 2498
$ws.Range("I132").Value = 2498
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7494
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4964
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2355.6
$ws.Range("J22").Value = 2694.25
$ws.Range("L22").Value = 2694.25
$ws.Range("N22").Value = -3284.25
$ws.Range("H27").Value = 2355.6
$ws.Range("J27").Value = 2694.25
$ws.Range("L27").Value = 2694.25
$ws.Range("N27").Value = -2908.25
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H46").Value = 4599.6665
$ws.Range("I46").Value = 5900
$ws.Range("K46").Value = 5900
$ws.Range("M46").Value = -5712
$ws.Range("H141").Value = 49715
$ws.Range("J141").Value = 49715
$ws.Range("L141").Value = 49715
$ws.Range("N141").Value = -60075

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
